# Update the Pulping Machine sheet with the revised values for the
# shock-related CO2 categories (ROI, Water Saving, Emission Saving, Land Saving)
# in row 4, as part of "Preparation for presentation of pulp mach shock".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G4").Value = -0.159999986179173
$ws.Range("H4").Value = -1
$ws.Range("I4").Value = -0.0007089301507221535
$ws.Range("J4").Value = -0.0003299660711491015
